# "Add new model" -- swap the Purchase-Order report's column headers and
# the single data row for a new set (new source model / run), and add a
# trailing "Payment Terms" style column (G) that didn't exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("A1").Value = "ORDER NO"
$ws.Range("B1").Value = "COMPANY NAME"
$ws.Range("C1").Value = "SUPPLIER NAME"
$ws.Range("D1").Value = "SUPPLIER NO"
$ws.Range("E1").Value = "PAGE TOTAL"
$ws.Range("F1").Value = "GRAND TOTAL"
$ws.Range("G1").Value = "PAYMENT TERMS"

# --- Data row (row 2) ---------------------------------------------------
# Column A is blank for this record.
$ws.Range("A2").ClearContents()

$ws.Range("B2").Value = "correctional services Department:"
$ws.Range("C2").Value = ":MOHLABANI CONSTRUCTION (PTY) LTD"
$ws.Range("D2").Value = "MAAA1327350"
$ws.Range("E2").Value = "R3 355 040.00"

# "3858" looks numeric, so force it to stay text (matches the other
# string-typed cells on this sheet) before writing it.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "3858"

$ws.Range("G2").Value = "STRICTLY WITHIN 3-5 WORKING DAYS"

# --- Cosmetic: widen the now 7-column table a bit ------------------------
$ws.Range("A1:G2").ColumnWidth = 20.78
